$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.154425621032715
$ws.Range("B1").Value = 2.381764650344849
$ws.Range("D1").Value = 2.395634889602661
$ws.Range("E1").Value = 1.225923538208008
